$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A163").Value = "06_22/23"
$ws.Range("B163").Value = 211.8
$ws.Range("C163").Value = 8

$ws.Range("A164").Value = "07_22/23"
$ws.Range("B164").Value = 234.4
$ws.Range("C164").Value = 9.2

$ws.Range("A165").Value = "08_22/23"
$ws.Range("B165").Value = 229.8
$ws.Range("C165").Value = 9.4
